$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions): remove cancelled event row 2 (昆山·ETHEREAL),
#     shift subsequent rows up, renumber the index column, and refresh the
#     "interested" counters (col F) that ticked up since the last scrape. ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("A2").EntireRow.Delete()
for ($r = 2; $r -le 39; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}
$ws1.Cells.Item(2, 6).Value = 604
$ws1.Cells.Item(5, 6).Value = 13
$ws1.Cells.Item(6, 6).Value = 15070
$ws1.Cells.Item(9, 6).Value = 682
$ws1.Cells.Item(10, 6).Value = 15253
$ws1.Cells.Item(11, 6).Value = 41
$ws1.Cells.Item(12, 6).Value = 8779
$ws1.Cells.Item(13, 6).Value = 341
$ws1.Cells.Item(15, 6).Value = 71
$ws1.Cells.Item(16, 6).Value = 184
$ws1.Cells.Item(19, 6).Value = 15
$ws1.Cells.Item(20, 6).Value = 27
$ws1.Cells.Item(21, 6).Value = 521
$ws1.Cells.Item(22, 6).Value = 23
$ws1.Cells.Item(24, 6).Value = 52
$ws1.Cells.Item(26, 6).Value = 9
$ws1.Cells.Item(28, 6).Value = 58
$ws1.Cells.Item(30, 6).Value = 32
$ws1.Cells.Item(34, 6).Value = 232
$ws1.Cells.Item(35, 6).Value = 279
$ws1.Cells.Item(38, 6).Value = 5393

# --- Sheet "演出" (shows): the 昆山·星月流光 interest counter ticked up too. ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 1006

# --- Sheet "全部类型" (all types, merged view): same row-2 removal + shift + refresh. ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("A2").EntireRow.Delete()
for ($r = 2; $r -le 42; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}
$ws4.Cells.Item(2, 6).Value = 604
$ws4.Cells.Item(5, 6).Value = 13
$ws4.Cells.Item(6, 6).Value = 15070
$ws4.Cells.Item(9, 6).Value = 682
$ws4.Cells.Item(10, 6).Value = 15253
$ws4.Cells.Item(11, 6).Value = 41
$ws4.Cells.Item(12, 6).Value = 8779
$ws4.Cells.Item(13, 6).Value = 341
$ws4.Cells.Item(15, 6).Value = 1006
$ws4.Cells.Item(16, 6).Value = 71
$ws4.Cells.Item(17, 6).Value = 184
$ws4.Cells.Item(20, 6).Value = 15
$ws4.Cells.Item(21, 6).Value = 27
$ws4.Cells.Item(22, 6).Value = 521
$ws4.Cells.Item(23, 6).Value = 23
$ws4.Cells.Item(25, 6).Value = 52
$ws4.Cells.Item(27, 6).Value = 9
$ws4.Cells.Item(29, 6).Value = 58
$ws4.Cells.Item(31, 6).Value = 32
$ws4.Cells.Item(37, 6).Value = 232
$ws4.Cells.Item(38, 6).Value = 279
$ws4.Cells.Item(41, 6).Value = 5393
